$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'23.259.57"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.64%  '
$ws.Range("D3").Value = "'1.615.69"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.71%  '
$ws.Range("D4").Value = "'1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.13%  '
$ws.Range("D5").Value = "'1.001"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.11%  '
$ws.Range("D6").Value = "'303.16"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.23%  '
$ws.Range("D7").Value = "'0.3781"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.17%  '
$ws.Range("D8").Value = "'51.60"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.94%  '
$ws.Range("D9").Value = "'0.3533"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.75%  '
$ws.Range("D10").Value = "'0.08112"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.16%  '
$ws.Range("D11").Value = "'1.207"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.61%  '
$ws.Range("D12").Value = "'1.001"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.12%  '
$ws.Range("D13").Value = "'22.21"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.94%  '
$ws.Range("D14").Value = "'6.363"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.86%  '
$ws.Range("D15").Value = "'7.277"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.75%  '
$ws.Range("D16").Value = "'0.00001212"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.66%  '
$ws.Range("D17").Value = "'1.605.35"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.88%  '
$ws.Range("D18").Value = "'94.12"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.59%  '
$ws.Range("D19").Value = "'0.06913"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.11%  '
$ws.Range("D20").Value = "'6.479"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.08%  '
$ws.Range("D21").Value = "'17.22"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.80%  '
$ws.Range("D22").Value = "'1.001"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.16%  '
$ws.Range("D23").Value = "'12.34"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.02%  '
$ws.Range("D24").Value = "'23.239.35"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.70%  '
$ws.Range("D25").Value = "'2.513"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.19%  '
$ws.Range("D26").Value = "'3.015"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -6.47%  '
$ws.Range("D27").Value = "'20.87"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.26%  '
$ws.Range("D28").Value = "'151.16"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.41%  '
$ws.Range("D29").Value = "'5.244"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.92%  '
$ws.Range("D30").Value = "'132.14"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.68%  '
$ws.Range("D31").Value = "'1.780.42"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.07%  '
$ws.Range("D32").Value = "'1.066"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +11.94%  '
$ws.Range("D33").Value = "'6.472"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -4.76%  '
$ws.Range("D34").Value = "'2.094"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -9.05%  '
$ws.Range("D35").Value = "'11.34"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.80%  '
$ws.Range("D36").Value = "'0.02708"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.15%  '
$ws.Range("D37").Value = "'0.08700"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.49%  '
$ws.Range("D38").Value = "'0.2454"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.05%  '
$ws.Range("D39").Value = "'0.06941"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.61%  '
$ws.Range("D40").Value = "'5.854"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -4.10%  '
$ws.Range("D41").Value = "'1.324"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.37%  '
$ws.Range("D42").Value = "'0.6884"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.66%  '
$ws.Range("D43").Value = "'11.95"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.02%  '
$ws.Range("D44").Value = "'15.29"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -5.96%  '
$ws.Range("D45").Value = "'1.000"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.11%  '
$ws.Range("D46").Value = "'0.6319"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.33%  '
$ws.Range("D47").Value = "'3.947"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.17%  '
$ws.Range("D48").Value = "'2.252"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.22%  '
$ws.Range("D49").Value = "'0.07868"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.57%  '
$ws.Range("D50").Value = "'127.41"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.41%  '
$ws.Range("D51").Value = "'1.170"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.89%  '
